$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D currently holds "Notes". Insert 3 new columns before it, pushing
# "Notes" from D to G, then fill the newly-opened D:F with the new headers.
$ws.Range("D1:F1").EntireColumn.Insert()

$ws.Range("D1").Value = "PartRevision"
$ws.Range("E1").Value = "ProcessRevision"
$ws.Range("F1").Value = "CustomerName"

# Match the formatting (style, width) of the existing header column.
$ws.Range("D1:F1").Style = $ws.Range("C1").Style
$ws.Range("D:F").ColumnWidth = $ws.Range("C:C").ColumnWidth
